$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (match formatting of the existing header cells, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("I1:J1").VerticalAlignment = -4160    # xlTop
$ws.Range("I1:J1").Borders.LineStyle = 1        # xlContinuous (thin box border)

# Data rows
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
